$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The source bank export pads amount text with non-breaking spaces (U+00A0),
# not regular spaces -- matters because TRIM() in the H-column formula only
# strips plain spaces, and the cached formula results below keep them.
$nbsp = [string][char]0x00A0
$pad2 = $nbsp + $nbsp

# Insert 5 new rows at the top (rows 1-5), shifting all existing data down.
$insertRange = $ws.Range("A1:A5")
$insertRange.EntireRow.Insert()

# Carry over the number formats (date / text styles) from the row that used
# to be row 1 (now row 6) onto the freshly inserted rows, so the new rows
# reuse the existing style indices instead of minting new ones.
$ws.Range("A6:G6").Copy()
$ws.Range("A1:G5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate columns A-F first, top row to bottom row (matches how the
# statement rows were keyed in) -- the 'saldo' (running balance) column G
# is filled in afterwards, oldest transaction first (see below), which is
# why the shared-string table ends up with the running-balance figures
# appended after all the concepto/documento/monto text.
$ws.Range("A1").Value = 41703
$ws.Range("B1").Value = "CONSUMO VISA NA KFC K104"
$ws.Range("C1").Value = "D"
$ws.Range("D1").Value = "0004085602"
$ws.Range("E1").Value = "INSTITUCIONAL SS.CC."
$ws.Range("F1").Value = "10.60" + $pad2

$ws.Range("A2").Value = 41703
$ws.Range("B2").Value = "CONSUMO DATA AKI MOLINEROS 161"
$ws.Range("C2").Value = "D"
$ws.Range("D2").Value = "0004018676"
$ws.Range("E2").Value = "INSTITUCIONAL SS.CC."
$ws.Range("F2").Value = "83.32" + $pad2

$ws.Range("A3").Value = 41698
$ws.Range("B3").Value = "INTERES A SU FAVOR"
$ws.Range("C3").Value = "C"
$ws.Range("D3").Value = "0000948985"
$ws.Range("E3").Value = "AGENCIA PARA PROCESOS BATCH"
$ws.Range("F3").Value = "1.53" + $pad2

$ws.Range("A4").Value = 41691
$ws.Range("B4").Value = "DEP CNB-1501119901001"
$ws.Range("C4").Value = "C"
$ws.Range("D4").Value = "0004766633"
$ws.Range("E4").Value = "AG. NORTE"
$ws.Range("F4").Value = "50.00" + $pad2

$ws.Range("A5").Value = 41691
$ws.Range("B5").Value = "DEP CNB-1501119901001"
$ws.Range("C5").Value = "C"
$ws.Range("D5").Value = "0004751756"
$ws.Range("E5").Value = "AG. NORTE"
$ws.Range("F5").Value = "50.00" + $pad2

# Running balance (saldo), entered oldest-transaction-first: row 5 up to row 1.
$ws.Range("G5").Value = "1943.95"
$ws.Range("G4").Value = "1993.95"
$ws.Range("G3").Value = "1995.48"
$ws.Range("G2").Value = "1912.16"
$ws.Range("G1").Value = "1901.56"

# The old H1/H2 formula cells shifted down to H6/H7 along with their rows;
# those rows no longer carry a formula in the updated layout, so clear them.
$ws.Range("H6:H7").ClearContents()

# Apply the shared formula across H1:H5 (it was previously only on H1:H2).
$ws.Range("H1:H5").Formula = "=CONCATENATE(" + [char]34 + "array('mo_fecha' => new \DateTime('" + [char]34 + ",TEXT(A1,""yyyy-mm-dd""),""'), 'mo_concepto' => '"",B1,""', 'mo_tipo' => '"",C1,""', 'mo_documento' => '"",D1,""', 'mo_oficina' => '"",E1,""', 'mo_monto' => "",TRIM(F1),"", 'mo_saldo' => "",G1,"", 'mo_fecha_crea' => new \DateTime('"",TEXT(NOW(),""yyyy-mm-dd H:m:s""),""'), 'mo_quien_crea' => 1, 'mo_fecha_modifica' => NULL, 'mo_quien_modifica' => NULL, 'mo_borrado_logico' => false),"")"

# Adjust selection/view as captured in the diff
$ws.Range("H1:H5").Select()
